$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database - new values for the two Austria utility rows (2 and 3)
foreach ($row in 2,3) {
    $ws.Range("E$row").Value = 0.0459
    $ws.Range("K$row").Value = 11.3
    $ws.Range("M$row").Value = 10.71
    $ws.Range("N$row").Value = 0.03536988110964333
    $ws.Range("O$row").Value = 0.9477876106194689
    $ws.Range("P$row").Value = 10.71
    $ws.Range("Q$row").Value = 0.03536988110964333
    $ws.Range("R$row").Value = 0.9477876106194689
    $ws.Range("U$row").Value = 0.011
    $ws.Range("V$row").Value = 0.00003632760898282695
    $ws.Range("W$row").Value = 0.1258351893095769
    $ws.Range("X$row").Value = 0.03429489876294754
    $ws.Range("Y$row").Value = 0.09154029054662931
    $ws.Range("AA$row").Value = -0.002937972468887202
    $ws.Range("AB$row").Value = 0.03421827093924794
    $ws.Range("AC$row").Value = -0.03715624340813514
    $ws.Range("AD$row").Value = 1.35
    $ws.Range("AF$row").Value = 1.35
    $ws.Range("AG$row").Value = 1.339
    $ws.Range("AH$row").Value = 0.004438599375308236
    $ws.Range("AI$row").Value = 0.01502504173622705
    $ws.Range("AJ$row").Value = 0.004402592235786927
    $ws.Range("AK$row").Value = 0.01490444016518439
    $ws.Range("AM$row").Value = -11.5
    $ws.Range("AQ$row").Value = 0.02321739130434783
}
